# The deck currently ships two theme parts:
#   theme1.xml -> "Office Theme" / "Office" colour scheme (orphaned - only
#                 wired to the Notes Master)
#   theme2.xml -> "Integral" / "Red Violet" colour scheme (the theme that is
#                 actually applied to the slide master + the presentation,
#                 i.e. what every slide renders with)
#
# The authored change swaps the two themes' contents: the live design
# (theme2.xml, used by every slide) switches from the "Integral"/Red Violet
# palette to the stock "Office Theme"/Office palette. Apply that by pushing
# the 12 standard Office theme colours into the presentation's active theme
# colour scheme (PowerPoint exposes this uniformly off any Slide object as
# ThemeColorScheme, items 1-12 = dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
